$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ScanSheet")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Sheet2 ("Sheet2"): clear rows 2-3 data, keep only C2 formatted (quote-prefix style like C15) ---
$ws2.Range("A2:H3").ClearContents()
$ws2.Range("C15").Copy()
$ws2.Range("C2").PasteSpecial(-4122)
$ws2.Range("A6").Select()

# --- Sheet1 ("ScanSheet"): only the selection changes, from C15 to B13 ---
$ws1.Activate()
$ws1.Range("B13").Select()
